$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 17:52"

$countries = @(
  "Estados Unidos",
  "España",
  "Italia",
  "Alemania",
  "Francia",
  "China",
  "Iran",
  "Reino Unido",
  "Turquia",
  "Suiza",
  "Belgica",
  "Paises Bajos",
  "Canada",
  "Austria",
  "Portugal",
  "Brasil",
  "Corea del Sur",
  "Israel",
  "Suecia",
  "Rusia",
  "Australia",
  "Noruega",
  "Irlanda",
  "India",
  "Chile",
  "Dinamarca",
  "Chequia",
  "Polonia",
  "Rumania",
  "Pakistan",
  "Malasia",
  "Japon",
  "Filipinas",
  "Ecuador",
  "Luxemburgo",
  "Arabia Saudita",
  "Indonesia",
  "Peru",
  "Serbia",
  "Mexico",
  "Finlandia",
  "Tailandia",
  "Panama",
  "Emiratos Arabes Unidos",
  "Republica Dominicana",
  "Catar",
  "Grecia",
  "Sudafrica",
  "Argentina",
  "Islandia",
  "Colombia",
  "Singapur",
  "Ucrania",
  "Argelia",
  "Egipto",
  "Croacia",
  "Nueva Zelanda",
  "Estonia",
  "Marruecos",
  "Irak",
  "Eslovenia",
  "Moldavia",
  "Hong Kong",
  "Lituania",
  "Bielorrusia",
  "Armenia",
  "Hungria",
  "Barein",
  "Bosnia y Herzegovina",
  "Kuwait",
  "Azerbaiyan",
  "Crucero",
  "Kazajistan",
  "Camerun",
  "Republica de Macedonia",
  "Tunez",
  "Eslovaquia",
  "Bulgaria",
  "Letonia",
  "Libano",
  "Principado de Andorra",
  "Uzbekistan",
  "Republica de Chipre",
  "Costa Rica",
  "Afganistan",
  "Uruguay",
  "Cuba",
  "Albania",
  "Taiwan",
  "Oman",
  "Burkina Faso",
  "Reunion",
  "Jordania",
  "Costa de Marfil",
  "Honduras",
  "Malta",
  "Ghana",
  "San Marino",
  "Mauricio",
  "Estado de Palestina",
  "Niger",
  "Vietnam",
  "Montenegro",
  "Nigeria",
  "Senegal",
  "Kirguistan",
  "Georgia",
  "Bolivia",
  "Sri Lanka",
  "Islas Feroe",
  "Kenia",
  "Venezuela",
  "Mayotte",
  "Banglades",
  "Consejo Danes para los Refugiados",
  "Martinica",
  "Isla de Man",
  "Guadalupe",
  "Brunei",
  "Guinea",
  "Paraguay",
  "Camboya",
  "Gibraltar",
  "Trinidad yTobago",
  "Ruanda",
  "Republica de Yibuti",
  "Madagascar",
  "El Salvador",
  "Monaco",
  "Guatemala",
  "Liechtenstein",
  "Guayana Francesa",
  "Aruba",
  "Togo",
  "Barbados",
  "Jamaica",
  "Mali",
  "Uganda",
  "Etiopia",
  "Congo",
  "Macao",
  "Polinesia Francesa",
  "Islas Caimanes",
  "Puerto Rico",
  "Zambia",
  "Bermudas",
  "San Martin (Parte Holandesa)",
  "Guinea-Bisau",
  "Bahamas",
  "Guam",
  "San Martin (Parte Francesa)",
  "Eritrea",
  "Guyana",
  "Gabon",
  "Benin",
  "Haiti",
  "Tanzania",
  "Birmania",
  "Libia",
  "Siria",
  "Maldivas",
  "Nueva Caledonia",
  "Islas Virgenes de los Estados Unidos",
  "Guinea Ecuatorial",
  "Namibia",
  "Angola",
  "Antigua y Barbuda",
  "Fiyi",
  "Dominica",
  "Mongolia",
  "Laos",
  "Santa Lucia",
  "Sudan",
  "Liberia",
  "Curazao",
  "Granada",
  "San Cristobal y Nieves",
  "Seychelles",
  "Groenlandia",
  "Zimbabue",
  "Surinam",
  "Mozambique",
  "Republica del Chad",
  "Suazilandia",
  "Nepal",
  "Montserrat",
  "Republica de Africa Central",
  "Islas Turcas y Caicos",
  "Malaui",
  "Santa Sede",
  "Belice",
  "San Vicente y las Granadinas",
  "Somalia",
  "Cabo Verde",
  "Sierra Leona",
  "Botsuana",
  "Nicaragua",
  "San Bartolome",
  "Mauritania",
  "Butan",
  "Sahara Occidental",
  "Santo Tome y Principe",
  "Gambia",
  "Islas Virgenes Britanicas",
  "Anguila",
  "Burundi",
  "Papua Nueva Guinea",
  "Bonaire, San Eustaquio y Saba",
  "Islas Malvinas",
  "San Pedro y Miquelon",
  "Sudan del Sur",
  "Timor Oriental"
)

$data = @(
  @(377538,10534,19904,345850,9015,913,11784),
  @(140511,3836,43208,83406,7069,556,13897),
  @(132547,0,22837,93187,3898,0,16523),
  @(105519,2144,36081,67584,4895,44,1854),
  @(98010,0,17250,71849,7072,0,8911),
  @(81740,32,77167,1242,211,0,3331),
  @(62589,2089,27039,31678,3987,133,3872),
  @(55242,3634,135,48948,1559,786,6159),
  @(30217,0,1326,28242,1415,0,649),
  @(22242,585,8056,13375,391,46,811),
  @(22194,1380,4157,16002,1260,403,2035),
  @(19580,777,250,17229,1424,234,2101),
  @(17046,379,3794,12907,426,22,345),
  @(12592,295,4046,8303,243,23,243),
  @(12442,712,184,11913,271,34,345),
  @(12345,162,127,11637,296,17,581),
  @(10331,47,6694,3445,55,6,192),
  @(9006,102,683,8263,153,3,60),
  @(7693,487,205,6897,640,114,591),
  @(7497,1154,494,6945,8,11,58),
  @(5919,24,2547,3324,93,3,48),
  @(5869,4,32,5749,83,12,88),
  @(5364,0,25,5165,165,0,174),
  @(5172,394,382,4653,0,1,137),
  @(5116,301,898,4175,337,6,43),
  @(5071,390,1491,3377,127,16,203),
  @(4944,122,147,4710,86,9,87),
  @(4532,119,191,4230,50,4,111),
  @(4417,360,460,3760,274,21,197),
  @(4005,239,429,3521,28,2,55),
  @(3963,170,1321,2579,92,1,63),
  @(3906,0,592,3222,79,0,92),
  @(3764,104,84,3503,1,14,177),
  @(3747,0,100,3456,156,0,191),
  @(2970,127,500,2426,35,3,44),
  @(2795,190,615,2139,41,3,41),
  @(2738,247,204,2313,0,12,221),
  @(2561,0,997,1472,89,0,92),
  @(2447,247,118,2268,109,3,61),
  @(2439,296,633,1681,89,31,125),
  @(2308,132,300,1974,81,7,34),
  @(2258,38,824,1407,30,1,27),
  @(2100,0,14,2031,88,0,55),
  @(2076,0,167,1898,1,0,11),
  @(1956,128,33,1825,147,12,98),
  @(1832,0,131,1697,37,0,4),
  @(1832,77,269,1482,90,2,81),
  @(1749,63,95,1641,7,1,13),
  @(1628,0,338,1235,96,2,55),
  @(1586,24,559,1021,11,0,6),
  @(1579,0,88,1445,50,0,46),
  @(1481,106,377,1098,29,0,6),
  @(1462,143,28,1389,16,7,45),
  @(1423,0,90,1160,46,0,173),
  @(1322,0,259,978,0,0,85),
  @(1282,60,167,1097,35,2,18),
  @(1160,54,241,918,14,0,1),
  @(1149,41,69,1059,12,2,21),
  @(1141,21,88,970,1,3,83),
  @(1122,91,373,684,0,1,65),
  @(1059,38,102,921,30,6,36),
  @(1056,91,40,994,80,3,22),
  @(936,21,236,696,12,0,4),
  @(880,37,8,857,11,0,15),
  @(861,161,54,794,31,0,13),
  @(853,20,87,758,30,0,8),
  @(817,73,71,699,17,9,47),
  @(811,55,458,349,4,0,4),
  @(754,80,68,653,4,4,33),
  @(743,78,105,637,23,0,1),
  @(717,76,44,665,23,1,8),
  @(712,0,619,82,10,0,11),
  @(685,23,50,629,16,0,6),
  @(658,0,17,632,0,0,9),
  @(599,29,30,543,15,3,26),
  @(596,0,25,549,39,0,22),
  @(581,47,13,566,3,0,2),
  @(577,28,42,512,21,1,23),
  @(548,6,16,530,5,1,2),
  @(548,7,60,469,27,0,19),
  @(525,0,31,473,12,0,21),
  @(504,47,30,472,8,0,2),
  @(494,29,45,440,11,0,9),
  @(467,0,18,447,14,0,2),
  @(423,56,18,391,0,3,14),
  @(415,0,123,286,14,0,6),
  @(396,46,27,358,12,2,11),
  @(383,6,131,230,7,1,22),
  @(376,3,61,310,0,0,5),
  @(371,40,67,302,3,0,2),
  @(364,0,108,238,0,0,18),
  @(349,0,40,309,4,0,0),
  @(349,0,126,217,5,0,6),
  @(323,0,41,279,0,0,3),
  @(305,7,6,277,10,0,22),
  @(293,52,5,288,4,0,0),
  @(287,73,31,251,2,0,5),
  @(277,0,35,210,14,0,32),
  @(268,24,8,253,3,0,7),
  @(261,7,42,218,0,0,1),
  @(253,0,26,217,0,0,10),
  @(249,4,123,126,8,0,0),
  @(241,8,4,235,7,0,2),
  @(238,0,35,198,2,0,5),
  @(237,11,105,130,1,0,2),
  @(228,12,33,191,5,0,4),
  @(195,7,45,147,6,1,3),
  @(194,11,2,178,3,3,14),
  @(185,7,42,137,5,1,6),
  @(184,1,120,64,2,0,0),
  @(172,14,7,159,2,0,6),
  @(165,0,65,93,6,0,7),
  @(164,0,15,147,3,0,2),
  @(164,41,33,114,1,5,17),
  @(161,0,5,138,0,0,18),
  @(151,0,50,97,20,0,4),
  @(150,11,73,76,6,0,1),
  @(139,0,31,101,14,0,7),
  @(135,0,85,49,3,0,1),
  @(128,0,5,123,0,0,0),
  @(115,2,15,95,1,0,5),
  @(115,1,58,57,1,0,0),
  @(113,4,60,53,0,0,0),
  @(106,1,1,97,0,0,8),
  @(105,0,4,101,0,0,0),
  @(90,0,9,81,0,0,0),
  @(82,0,2,80,6,0,0),
  @(78,9,5,69,4,0,4),
  @(77,0,4,72,4,0,1),
  @(77,7,17,57,3,0,3),
  @(77,0,55,21,0,0,1),
  @(72,0,34,38,1,0,0),
  @(71,0,2,69,0,0,0),
  @(65,7,23,39,0,0,3),
  @(60,0,6,52,4,0,2),
  @(59,0,8,48,0,0,3),
  @(56,9,12,39,0,0,5),
  @(52,0,0,52,0,0,0),
  @(52,8,4,46,1,0,2),
  @(45,0,2,38,0,0,5),
  @(44,0,10,34,1,0,0),
  @(42,0,0,42,0,0,0),
  @(39,0,1,37,0,0,1),
  @(39,0,1,36,0,0,2),
  @(39,0,7,31,0,0,1),
  @(39,0,17,20,0,0,2),
  @(37,0,1,30,0,0,6),
  @(33,15,0,33,0,0,0),
  @(33,0,5,23,1,0,5),
  @(32,0,0,31,0,0,1),
  @(32,0,7,23,6,0,2),
  @(31,0,0,31,0,0,0),
  @(31,0,8,18,8,1,5),
  @(30,6,1,28,0,0,1),
  @(26,0,5,20,0,0,1),
  @(25,1,0,24,0,0,1),
  @(24,0,5,18,0,0,1),
  @(22,0,0,21,0,0,1),
  @(19,0,1,17,0,0,1),
  @(19,0,2,15,0,0,2),
  @(19,0,13,6,0,0,0),
  @(18,0,1,17,0,0,0),
  @(17,0,0,17,0,0,0),
  @(16,0,3,13,0,0,0),
  @(16,0,3,13,0,0,0),
  @(16,0,2,12,0,0,2),
  @(15,0,0,15,1,0,0),
  @(15,1,0,15,0,0,0),
  @(15,0,1,14,0,0,0),
  @(15,0,4,11,0,0,0),
  @(14,2,0,14,0,0,0),
  @(14,0,1,13,0,0,0),
  @(14,2,2,10,0,0,2),
  @(14,0,3,8,0,0,3),
  @(13,0,5,7,0,0,1),
  @(12,0,0,12,2,0,0),
  @(11,1,0,11,0,0,0),
  @(11,0,0,11,0,0,0),
  @(11,0,10,1,0,0,0),
  @(10,0,0,9,0,0,1),
  @(10,0,0,9,0,0,1),
  @(10,0,1,9,0,0,0),
  @(10,1,2,8,0,0,0),
  @(10,0,4,6,0,0,0),
  @(9,0,1,8,0,0,0),
  @(9,0,0,7,0,0,2),
  @(8,0,0,8,0,0,0),
  @(8,0,0,7,0,0,1),
  @(8,3,0,7,1,1,1),
  @(7,0,0,7,0,0,0),
  @(7,0,0,6,1,0,1),
  @(7,0,1,6,0,0,0),
  @(7,0,1,6,0,0,0),
  @(7,0,1,5,0,0,1),
  @(6,0,0,6,0,0,0),
  @(6,0,0,5,0,0,1),
  @(6,0,0,5,0,0,1),
  @(6,0,1,5,0,0,0),
  @(6,0,2,3,0,0,1),
  @(5,0,2,3,0,0,0),
  @(4,0,0,4,0,0,0),
  @(4,0,0,4,0,0,0),
  @(4,0,2,1,0,0,1),
  @(3,0,0,3,0,0,0),
  @(3,0,0,3,0,0,0),
  @(3,0,0,3,0,0,0),
  @(2,0,0,2,0,0,0),
  @(2,0,0,2,0,0,0),
  @(2,0,0,2,0,0,0),
  @(1,0,0,1,0,0,0),
  @(1,0,0,1,0,0,0),
  @(1,0,0,1,0,0,0)
)

for ($i = 0; $i -lt $countries.Length; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $countries[$i]
    $rowvals = $data[$i]
    for ($c = 0; $c -lt $rowvals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $rowvals[$c]
    }
}
